# =====================================================================
# Apply the "Add files via upload" commit to the workbook:
#  - fix the spelling of "litterature" -> "literature" (sheet name +
#    the _FilterDatabase defined name follow automatically)
#  - add a "Source" column to both data sheets, citing the paper/
#    dataset each metal's projection came from
#  - add a new "Ref" sheet with the full reference list + hyperlinks
#  - leave "OSD by scenario" as the active/selected sheet
# =====================================================================

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # OSD proj literature
$ws2 = $wb.Worksheets.Item(2)   # OSD by scenario

# ---------------------------------------------------------------------
# 1) Fix sheet-1's name (autoFilter defined name text updates with it)
# ---------------------------------------------------------------------
$ws1.Name = "OSD proj literature"

# ---------------------------------------------------------------------
# 2) "Source" header column on sheet 1 (col F) and sheet 2 (col G)
#    Formatting = same bold / light-blue header style already used for
#    the other headers, just with a left+right border only (no
#    top/bottom) since there is no box drawn under the Source column.
# ---------------------------------------------------------------------
function Style-SourceHeader($cell, $templateCell) {
    $templateCell.Copy()
    $cell.PasteSpecial(-4122)   # xlPasteFormats
    $excel.CutCopyMode = $false
    $cell.Borders.Item(8).LineStyle = -4142   # xlEdgeTop    -> none
    $cell.Borders.Item(9).LineStyle = -4142   # xlEdgeBottom -> none
}

$ws1.Range("F1").Value = "Source"
Style-SourceHeader $ws1.Range("F1") $ws1.Range("B1")

$ws2.Range("G1").Value = "Source"
Style-SourceHeader $ws2.Range("G1") $ws2.Range("F1")

# ---------------------------------------------------------------------
# 3) New "Ref" sheet, appended after "OSD by scenario". Build its
#    reference rows now (this is the order the strings were first
#    authored in, so the shared-string table matches the source file).
# ---------------------------------------------------------------------
$refSheet = $wb.Worksheets.Add($null, $ws2)
$refSheet.Name = "Ref"

$refSheet.Range("A2").Value = "(Pedneault et al., 2022)"
$refSheet.Range("B2").Value = "Pedneault, J., Majeau-Bettez, G., Pauliuk, S., & Margni, M. (2022). Sector-specific scenarios for future stocks and flows of aluminum : An analysis based on shared socioeconomic pathways. Journal of Industrial Ecology, 26(5), 1728-1746. https://doi.org/10.1111/jiec.13321"

$refSheet.Range("A3").Value = "(International Energy Agency, 2024)"
$refSheet.Range("B3").Value = "International Energy Agency. (2024). Critical minerals data explorer - Data tools. https://www.iea.org/data-and-statistics/data-tools/critical-minerals-data-explorer"

$refSheet.Range("A4").Value = "(Kermeli et al., 2022)"
$refSheet.Range("B4").Value = "Kermeli, K., Edelenbosch, O. Y., Crijns-Graus, W., Van Ruijven, B. J., Van Vuuren, D. P., & Worrell, E. (2022). Improving material projections in Integrated Assessment Models : The use of a stock-based versus a flow-based approach for the iron and steel industry. Energy, 239, 122434. https://doi.org/10.1016/j.energy.2021.122434"

$refSheet.Range("A5").Value = "(Rostek et al., 2023)"
$refSheet.Range("B5").Value = "Rostek, L., Pirard, E., & Loibl, A. (2023). The future availability of zinc : Potential contributions from recycling and necessary ones from mining. Resources, Conservation & Recycling Advances, 19, 200166. https://doi.org/10.1016/j.rcradv.2023.200166"

$refSheet.Range("A1").Value = "Ref"

# Title band A1:B1 ("Ref"), merged, centered, light-grey fill, boxed
$titleRng = $refSheet.Range("A1:B1")
$titleRng.Merge()
$titleRng.HorizontalAlignment = -4108      # xlCenter
$titleRng.Interior.Pattern = 1             # xlSolid
$titleRng.Interior.Color = 15921906        # light grey (~ theme0 tint -5%)
$titleRng.Borders.Item(7).LineStyle  = 1   # xlEdgeLeft
$titleRng.Borders.Item(10).LineStyle = 1   # xlEdgeRight
$titleRng.Borders.Item(8).LineStyle  = 1   # xlEdgeTop
$titleRng.Borders.Item(9).LineStyle  = 1   # xlEdgeBottom

# Hyperlinks for column B (display text kept as the shared long citation,
# the hyperlink target is the article / dataset URL)
$refSheet.Hyperlinks.Add($refSheet.Range("B2"), "https://doi.org/10.1111/jiec.13321") | Out-Null
$refSheet.Hyperlinks.Add($refSheet.Range("B3"), "https://www.iea.org/data-and-statistics/data-tools/critical-minerals-data-explorer") | Out-Null
$refSheet.Hyperlinks.Add($refSheet.Range("B4"), "https://doi.org/10.1016/j.energy.2021.122434") | Out-Null
$refSheet.Hyperlinks.Add($refSheet.Range("B5"), "https://doi.org/10.1016/j.rcradv.2023.200166") | Out-Null

# Indent the three single-line hyperlink cells (B2, B4, B5); B3 keeps the
# plain hyperlink style with no indent, matching the source formatting.
$refSheet.Range("B2").IndentLevel = 2
$refSheet.Range("B4").IndentLevel = 2
$refSheet.Range("B5").IndentLevel = 2
$refSheet.Range("B2").VerticalAlignment = -4108   # xlCenter
$refSheet.Range("B4").VerticalAlignment = -4108
$refSheet.Range("B5").VerticalAlignment = -4108

$refSheet.Columns.Item(2).ColumnWidth = 17.26953125

# ---------------------------------------------------------------------
# 4) Per-row source citations (reuse the strings already registered
#    above, so no new shared-string entries get created here)
#    Sheet 1 "OSD proj literature": Iron -> Kermeli et al., Zinc -> Rostek et al.
# ---------------------------------------------------------------------
$ws1.Range("F2").Value = "(Kermeli et al., 2022)"
$ws1.Range("F3").Value = "(Rostek et al., 2023)"

#    Sheet 2 "OSD by scenario": Aluminium rows (2-5) -> Pedneault et al.
#    Cobalt / Copper / Nickel rows (6-14) -> International Energy Agency
for ($r = 2; $r -le 5; $r++) {
    $ws2.Range("G$r").Value = "(Pedneault et al., 2022)"
}
for ($r = 6; $r -le 14; $r++) {
    $ws2.Range("G$r").Value = "(International Energy Agency, 2024)"
}

# ---------------------------------------------------------------------
# 5) Final selection / active-sheet state: "OSD by scenario" tab active
# ---------------------------------------------------------------------
$ws1.Range("F3").Select()
$refSheet.Range("C12").Select()
$ws2.Activate()
$ws2.Range("E18").Select()
